# Updated symbol list on Fri Dec 23 22:22:45 UTC 2022 with GitHub Actions
#
# Re-applies the latest crypto price/volume scrape onto Sheet1. All data
# cells in this sheet are stored as text (t="inlineStr"/shared-string),
# even the numeric-looking "Price" column, so every write below goes
# through a small helper that forces text entry (leading apostrophe) and
# then resets the cell style back to Normal so no stray NumberFormat /
# quotePrefix style sticks around on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Cell, $Text) {
    $range = $Sheet.Range($Cell)
    # Leading apostrophe forces Excel to treat the entry as literal text
    # instead of re-parsing numeric-looking strings (e.g. "246.02") into
    # a number.
    $range.Value = "'" + $Text
    # The apostrophe entry stamps a quotePrefix style on the cell; put the
    # cell style back to Normal so the underlying style index is unchanged.
    $range.Style = "Normal"
}

# ---- Price column (D) updates -------------------------------------------------
Set-TextValue $ws "D2"  "246.02"
Set-TextValue $ws "D3"  "22.10"
Set-TextValue $ws "D4"  "5.326"
Set-TextValue $ws "D5"  "0.05870"
Set-TextValue $ws "D6"  "3.384"
Set-TextValue $ws "D7"  "6.384"
Set-TextValue $ws "D8"  "0.8171"
Set-TextValue $ws "D9"  "0.9530"
Set-TextValue $ws "D11" "0.03519"
Set-TextValue $ws "D12" "0.07348"
Set-TextValue $ws "D13" "0.03039"
Set-TextValue $ws "D14" "4.427"
Set-TextValue $ws "D15" "0.09399"
Set-TextValue $ws "D16" "0.001585"
Set-TextValue $ws "D17" "0.04808"
Set-TextValue $ws "D18" "0.0005904"
Set-TextValue $ws "D19" "0.006050"
Set-TextValue $ws "D20" "0.004085"
Set-TextValue $ws "D21" "0.0009858"
Set-TextValue $ws "D22" "0.00009710"
Set-TextValue $ws "D23" "3.683"
Set-TextValue $ws "D24" "2.179"
Set-TextValue $ws "D25" "0.3256"
Set-TextValue $ws "D26" "0.1261"
Set-TextValue $ws "D27" "0.0002473"
Set-TextValue $ws "D40" "0.03868"

# ---- Rows 41-43: the coin ranking rotated by one position ---------------------
# old row41 (BKEXToken) data -> now at row42
# old row42 (CEJI)       data -> now at row43
# old row43 (KickToken)  data -> now at row41 (and loses its "Worstin24h" flag)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006617"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1075"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002442"
$ws.Range("E43").Value = "42CEJICEJI"

# ---- Remaining price / flag updates -------------------------------------------
Set-TextValue $ws "D44" "0.005317"
Set-TextValue $ws "D45" "0.00005676"
Set-TextValue $ws "D47" "0.7755"
Set-TextValue $ws "D48" "0.07111"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
Set-TextValue $ws "D49" "0.00002101"
Set-TextValue $ws "D50" "0.01011"
